$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 38.92704733333333
$ws.Range("H2").Value = 116.781142
$ws.Range("I2").Value = 0.2481365284058833
$ws.Range("J2").Value = 0.2481365284058833
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 0.6763629999999999
$ws.Range("N2").Value = 2.029089
$ws.Range("O2").Value = 0.6127318215515719
$ws.Range("P2").Value = 0.6127318215515719
$ws.Range("Q2").Value = 26.32881451551533
$ws.Range("R2").Value = 236.959330639638
$ws.Range("S2").Value = 0.1520411470436202
$ws.Range("T2").Value = 0.1520411470436202

$ws.Range("G3").Value = 38.92704733333333
$ws.Range("H3").Value = 116.781142
$ws.Range("I3").Value = 0.2481365284058833
$ws.Range("J3").Value = 0.2481365284058833
$ws.Range("M3").Value = 0.05377866666666667
$ws.Range("O3").Value = 0.04871925339984812
$ws.Range("P3").Value = 0.04871925339984811
$ws.Range("Q3").Value = 2.093444702856889
$ws.Range("R3").Value = 18.841002325712
$ws.Range("S3").Value = 0.01208902640516484
$ws.Range("T3").Value = 0.01208902640516484

$ws.Range("G4").Value = 38.92704733333333
$ws.Range("H4").Value = 116.781142
$ws.Range("I4").Value = 0.2481365284058833
$ws.Range("J4").Value = 0.2481365284058833
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 0.3737066666666666
$ws.Range("N4").Value = 1.12112
$ws.Range("O4").Value = 0.3385489250485801
$ws.Range("P4").Value = 0.33854892504858
$ws.Range("Q4").Value = 14.54729710211555
$ws.Range("R4").Value = 130.92567391904
$ws.Range("S4").Value = 0.08400635495709825
$ws.Range("T4").Value = 0.08400635495709824

$ws.Range("G5").Value = 97.97927366666669
$ws.Range("I5").Value = 0.6245589760556541
$ws.Range("J5").Value = 0.6245589760556541
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 0.6763629999999999
$ws.Range("N5").Value = 2.029089
$ws.Range("O5").Value = 0.6127318215515719
$ws.Range("P5").Value = 0.6127318215515719
$ws.Range("Q5").Value = 66.26955547500768
$ws.Range("R5").Value = 596.4259992750691
$ws.Range("S5").Value = 0.3826871590649655
$ws.Range("T5").Value = 0.3826871590649655

$ws.Range("G6").Value = 97.97927366666669
$ws.Range("I6").Value = 0.6245589760556541
$ws.Range("J6").Value = 0.6245589760556541
$ws.Range("M6").Value = 0.05377866666666667
$ws.Range("O6").Value = 0.04871925339984812
$ws.Range("P6").Value = 0.04871925339984811
$ws.Range("Q6").Value = 5.269194698761779
$ws.Range("R6").Value = 47.42275228885601
$ws.Range("S6").Value = 0.03042804701760509
$ws.Range("T6").Value = 0.03042804701760508

$ws.Range("G7").Value = 97.97927366666669
$ws.Range("I7").Value = 0.6245589760556541
$ws.Range("J7").Value = 0.6245589760556541
$ws.Range("K7").Value = 3
$ws.Range("L7").Value = 1
$ws.Range("M7").Value = 0.3737066666666666
$ws.Range("N7").Value = 1.12112
$ws.Range("O7").Value = 0.3385489250485801
$ws.Range("P7").Value = 0.33854892504858
$ws.Range("Q7").Value = 36.61550776439111
$ws.Range("R7").Value = 329.53956987952
$ws.Range("S7").Value = 0.2114437699730836
$ws.Range("T7").Value = 0.2114437699730835

$ws.Range("G8").Value = 19.463074
$ws.Range("H8").Value = 58.389222
$ws.Range("I8").Value = 0.1240653978482281
$ws.Range("J8").Value = 0.1240653978482281
$ws.Range("K8").Value = 3
$ws.Range("L8").Value = 1
$ws.Range("M8").Value = 0.6763629999999999
$ws.Range("N8").Value = 2.029089
$ws.Range("O8").Value = 0.6127318215515719
$ws.Range("P8").Value = 0.6127318215515719
$ws.Range("Q8").Value = 13.164103119862
$ws.Range("R8").Value = 118.476928078758
$ws.Range("S8").Value = 0.0760188172150653
$ws.Range("T8").Value = 0.0760188172150653

$ws.Range("G9").Value = 19.463074
$ws.Range("H9").Value = 58.389222
$ws.Range("I9").Value = 0.1240653978482281
$ws.Range("J9").Value = 0.1240653978482281
$ws.Range("M9").Value = 0.05377866666666667
$ws.Range("O9").Value = 0.04871925339984812
$ws.Range("P9").Value = 0.04871925339984811
$ws.Range("Q9").Value = 1.046698168954667
$ws.Range("R9").Value = 9.420283520592001
$ws.Range("S9").Value = 0.006044373555920799
$ws.Range("T9").Value = 0.006044373555920798

$ws.Range("G10").Value = 19.463074
$ws.Range("H10").Value = 58.389222
$ws.Range("I10").Value = 0.1240653978482281
$ws.Range("J10").Value = 0.1240653978482281
$ws.Range("K10").Value = 3
$ws.Range("L10").Value = 1
$ws.Range("M10").Value = 0.3737066666666666
$ws.Range("N10").Value = 1.12112
$ws.Range("O10").Value = 0.3385489250485801
$ws.Range("P10").Value = 0.33854892504858
$ws.Range("Q10").Value = 7.273480507626667
$ws.Range("R10").Value = 65.46132456863999
$ws.Range("S10").Value = 0.04200220707724206
$ws.Range("T10").Value = 0.04200220707724205

$ws.Range("G11").Value = 0.5081416666666667
$ws.Range("H11").Value = 1.524425
$ws.Range("I11").Value = 0.003239097690234427
$ws.Range("J11").Value = 0.003239097690234427
$ws.Range("K11").Value = 3
$ws.Range("L11").Value = 1
$ws.Range("M11").Value = 0.6763629999999999
$ws.Range("N11").Value = 2.029089
$ws.Range("O11").Value = 0.6127318215515719
$ws.Range("P11").Value = 0.6127318215515719
$ws.Range("Q11").Value = 0.3436882220916667
$ws.Range("R11").Value = 3.093193998825
$ws.Range("S11").Value = 0.00198469822792083
$ws.Range("T11").Value = 0.00198469822792083

$ws.Range("G12").Value = 0.5081416666666667
$ws.Range("H12").Value = 1.524425
$ws.Range("I12").Value = 0.003239097690234427
$ws.Range("J12").Value = 0.003239097690234427
$ws.Range("M12").Value = 0.05377866666666667
$ws.Range("O12").Value = 0.04871925339984812
$ws.Range("P12").Value = 0.04871925339984811
$ws.Range("Q12").Value = 0.02732718131111112
$ws.Range("R12").Value = 0.2459446318
$ws.Range("S12").Value = 0.0001578064211573938
$ws.Range("T12").Value = 0.0001578064211573938

$ws.Range("G13").Value = 0.5081416666666667
$ws.Range("H13").Value = 1.524425
$ws.Range("I13").Value = 0.003239097690234427
$ws.Range("J13").Value = 0.003239097690234427
$ws.Range("K13").Value = 3
$ws.Range("L13").Value = 1
$ws.Range("M13").Value = 0.3737066666666666
$ws.Range("N13").Value = 1.12112
$ws.Range("O13").Value = 0.3385489250485801
$ws.Range("P13").Value = 0.33854892504858
$ws.Range("Q13").Value = 0.1898959284444444
$ws.Range("R13").Value = 1.709063356
$ws.Range("S13").Value = 0.001096593041156204
$ws.Range("T13").Value = 0.001096593041156204

